# Insert a new weekly price record for Cilantro (Agrícola del Norte S.A. de Arica)
# right above the current row 31, pushing the existing rows 31-48 down to 32-49.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 31 (shifts rows 31..48 -> 32..49,
# carrying along formatting such as the date style on column D).
$ws.Rows.Item(31).Insert()

# Populate the newly inserted row 31 with the new observation.
$ws.Cells.Item(31, 1).Value = 1
$ws.Cells.Item(31, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(31, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(31, 4).Value = 44455
$ws.Cells.Item(31, 5).Value = 15
$ws.Cells.Item(31, 6).Value = 100112040
$ws.Cells.Item(31, 7).Value = "Cilantro"
$ws.Cells.Item(31, 8).Value = "Sin especificar"
$ws.Cells.Item(31, 9).Value = "Primera"
$ws.Cells.Item(31, 10).Value = 250
$ws.Cells.Item(31, 11).Value = 900
$ws.Cells.Item(31, 12).Value = 1000
$ws.Cells.Item(31, 13).Value = 950
$ws.Cells.Item(31, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(31, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(31, 16).Value = 475
$ws.Cells.Item(31, 17).Value = 2
$ws.Cells.Item(31, 18).Value = "Hortaliza"
